$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 10 and 12 swapped their full match data (id column A stays put) ---
# Row 10
$ws.Cells.Item(10, 2).Value = 6627290
$ws.Cells.Item(10, 3).Value = "Bulgaria First League"
$ws.Cells.Item(10, 4).Value = "Bulgaria First League"
$ws.Cells.Item(10, 5).Value = 45084.53125
$ws.Cells.Item(10, 6).Value = "Cherno More Varna"
$ws.Cells.Item(10, 7).Value = "Ludogorets Razgrad"
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 9).Value = 1
$ws.Cells.Item(10, 10).Value = "A"
$ws.Cells.Item(10, 11).Value = 6
$ws.Cells.Item(10, 12).Value = 4
$ws.Cells.Item(10, 13).Value = 1.45
$ws.Cells.Item(10, 14).Value = 6
$ws.Cells.Item(10, 15).Value = 4.333
$ws.Cells.Item(10, 16).Value = 1.55
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 19).Value = 1.85
$ws.Cells.Item(10, 20).Value = 2.75
$ws.Cells.Item(10, 21).Value = 1.825
$ws.Cells.Item(10, 22).Value = 2.025
$ws.Cells.Item(10, 23).Value = -1
$ws.Cells.Item(10, 24).Value = -1
$ws.Cells.Item(10, 25).Value = 0.55
$ws.Cells.Item(10, 26).Value = 0
$ws.Cells.Item(10, 27).Value = -0
$ws.Cells.Item(10, 28).Value = -1
$ws.Cells.Item(10, 29).Value = 1.025

# Row 12
$ws.Cells.Item(12, 2).Value = 6627724
$ws.Cells.Item(12, 3).Value = "Bulgaria First League"
$ws.Cells.Item(12, 4).Value = "Bulgaria First League"
$ws.Cells.Item(12, 5).Value = 45084.53125
$ws.Cells.Item(12, 6).Value = "CSKA 1948 Sofia"
$ws.Cells.Item(12, 7).Value = "Lokomotiv Plovdiv"
$ws.Cells.Item(12, 8).Value = 0
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = "D"
$ws.Cells.Item(12, 11).Value = 1.5
$ws.Cells.Item(12, 12).Value = 3.8
$ws.Cells.Item(12, 13).Value = 6
$ws.Cells.Item(12, 14).Value = 1.45
$ws.Cells.Item(12, 15).Value = 4.2
$ws.Cells.Item(12, 16).Value = 8
$ws.Cells.Item(12, 17).Value = -1.25
$ws.Cells.Item(12, 18).Value = 2.025
$ws.Cells.Item(12, 19).Value = 1.825
$ws.Cells.Item(12, 20).Value = 2.5
$ws.Cells.Item(12, 21).Value = 1.85
$ws.Cells.Item(12, 22).Value = 2
$ws.Cells.Item(12, 23).Value = -1
$ws.Cells.Item(12, 24).Value = 3.2
$ws.Cells.Item(12, 25).Value = -1
$ws.Cells.Item(12, 26).Value = -1
$ws.Cells.Item(12, 27).Value = 0.825
$ws.Cells.Item(12, 28).Value = -1
$ws.Cells.Item(12, 29).Value = 1

# --- Rows 241 and 242: existing rows, data updated in place (values only; formatting already present) ---
# Row 241
$ws.Cells.Item(241, 2).Value = 6978454
$ws.Cells.Item(241, 3).Value = "Bulgaria First League"
$ws.Cells.Item(241, 4).Value = "Bulgaria First League"
$ws.Cells.Item(241, 5).Value = 45395.59375
$ws.Cells.Item(241, 6).Value = "Cherno More Varna"
$ws.Cells.Item(241, 7).Value = "Levski Sofia"
$ws.Cells.Item(241, 8).Value = 3
$ws.Cells.Item(241, 9).Value = 1
$ws.Cells.Item(241, 10).Value = "H"
$ws.Cells.Item(241, 11).Value = 2.7
$ws.Cells.Item(241, 12).Value = 3.1
$ws.Cells.Item(241, 13).Value = 2.7
$ws.Cells.Item(241, 14).Value = 2.3
$ws.Cells.Item(241, 15).Value = 3.1
$ws.Cells.Item(241, 16).Value = 3.4
$ws.Cells.Item(241, 17).Value = -0.25
$ws.Cells.Item(241, 18).Value = 1.95
$ws.Cells.Item(241, 19).Value = 1.9
$ws.Cells.Item(241, 20).Value = 2
$ws.Cells.Item(241, 21).Value = 1.875
$ws.Cells.Item(241, 22).Value = 1.975
$ws.Cells.Item(241, 23).Value = 1.3
$ws.Cells.Item(241, 24).Value = -1
$ws.Cells.Item(241, 25).Value = -1
$ws.Cells.Item(241, 26).Value = 0.95
$ws.Cells.Item(241, 27).Value = -1
$ws.Cells.Item(241, 28).Value = 0.875
$ws.Cells.Item(241, 29).Value = -1

# Row 242
$ws.Cells.Item(242, 2).Value = 6978457
$ws.Cells.Item(242, 3).Value = "Bulgaria First League"
$ws.Cells.Item(242, 4).Value = "Bulgaria First League"
$ws.Cells.Item(242, 5).Value = 45396.45833333334
$ws.Cells.Item(242, 6).Value = "Etar 1924 Veliko Tarnovo"
$ws.Cells.Item(242, 7).Value = "Krumovgrad"
$ws.Cells.Item(242, 8).Value = 0
$ws.Cells.Item(242, 9).Value = 0
$ws.Cells.Item(242, 10).Value = "D"
$ws.Cells.Item(242, 11).Value = 4
$ws.Cells.Item(242, 12).Value = 3.2
$ws.Cells.Item(242, 13).Value = 2
$ws.Cells.Item(242, 14).Value = 5.5
$ws.Cells.Item(242, 15).Value = 3.6
$ws.Cells.Item(242, 16).Value = 1.7
$ws.Cells.Item(242, 17).Value = 0.75
$ws.Cells.Item(242, 18).Value = 2
$ws.Cells.Item(242, 19).Value = 1.85
$ws.Cells.Item(242, 20).Value = 2.25
$ws.Cells.Item(242, 21).Value = 2
$ws.Cells.Item(242, 22).Value = 1.85
$ws.Cells.Item(242, 23).Value = -1
$ws.Cells.Item(242, 24).Value = 2.6
$ws.Cells.Item(242, 25).Value = -1
$ws.Cells.Item(242, 26).Value = 1
$ws.Cells.Item(242, 27).Value = -1
$ws.Cells.Item(242, 28).Value = -1
$ws.Cells.Item(242, 29).Value = 0.8500000000000001

# --- Rows 243-247: brand-new rows appended; copy formatting from row 240 first, then set values ---
# Row 243
$ws.Range("A240:AC240").Copy() | Out-Null
$ws.Range("A243:AC243").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(243, 1).Value = 241
$ws.Cells.Item(243, 2).Value = 6978458
$ws.Cells.Item(243, 3).Value = "Bulgaria First League"
$ws.Cells.Item(243, 4).Value = "Bulgaria First League"
$ws.Cells.Item(243, 5).Value = 45401.47916666666
$ws.Cells.Item(243, 6).Value = "Botev Vratsa"
$ws.Cells.Item(243, 7).Value = "Etar 1924 Veliko Tarnovo"
$ws.Cells.Item(243, 11).Value = 1.8
$ws.Cells.Item(243, 12).Value = 3.3
$ws.Cells.Item(243, 13).Value = 4.75
$ws.Cells.Item(243, 14).Value = 1.7
$ws.Cells.Item(243, 15).Value = 3.3
$ws.Cells.Item(243, 16).Value = 5.5
$ws.Cells.Item(243, 17).Value = -0.75
$ws.Cells.Item(243, 18).Value = 1.925
$ws.Cells.Item(243, 19).Value = 1.925
$ws.Cells.Item(243, 20).Value = 2.25
$ws.Cells.Item(243, 21).Value = 1.825
$ws.Cells.Item(243, 22).Value = 2.025
$ws.Cells.Item(243, 23).Value = 0
$ws.Cells.Item(243, 24).Value = 0
$ws.Cells.Item(243, 25).Value = 0
$ws.Cells.Item(243, 26).Value = 0
$ws.Cells.Item(243, 27).Value = 0

# Row 244
$ws.Range("A240:AC240").Copy() | Out-Null
$ws.Range("A244:AC244").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(244, 1).Value = 242
$ws.Cells.Item(244, 2).Value = 6978459
$ws.Cells.Item(244, 3).Value = "Bulgaria First League"
$ws.Cells.Item(244, 4).Value = "Bulgaria First League"
$ws.Cells.Item(244, 5).Value = 45401.58333333334
$ws.Cells.Item(244, 6).Value = "Krumovgrad"
$ws.Cells.Item(244, 7).Value = "Slavia Sofia"
$ws.Cells.Item(244, 11).Value = 2.3
$ws.Cells.Item(244, 12).Value = 3.1
$ws.Cells.Item(244, 13).Value = 3.25
$ws.Cells.Item(244, 14).Value = 2.8
$ws.Cells.Item(244, 15).Value = 3
$ws.Cells.Item(244, 16).Value = 2.625
$ws.Cells.Item(244, 17).Value = 0
$ws.Cells.Item(244, 18).Value = 2
$ws.Cells.Item(244, 19).Value = 1.85
$ws.Cells.Item(244, 20).Value = 2
$ws.Cells.Item(244, 21).Value = 1.775
$ws.Cells.Item(244, 22).Value = 2.1
$ws.Cells.Item(244, 23).Value = 0
$ws.Cells.Item(244, 24).Value = 0
$ws.Cells.Item(244, 25).Value = 0
$ws.Cells.Item(244, 26).Value = 0
$ws.Cells.Item(244, 27).Value = 0

# Row 245
$ws.Range("A240:AC240").Copy() | Out-Null
$ws.Range("A245:AC245").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(245, 1).Value = 243
$ws.Cells.Item(245, 2).Value = 6978461
$ws.Cells.Item(245, 3).Value = "Bulgaria First League"
$ws.Cells.Item(245, 4).Value = "Bulgaria First League"
$ws.Cells.Item(245, 5).Value = 45402.38541666666
$ws.Cells.Item(245, 6).Value = "Pirin Blagoevgrad"
$ws.Cells.Item(245, 7).Value = "Cherno More Varna"
$ws.Cells.Item(245, 11).Value = 5.75
$ws.Cells.Item(245, 12).Value = 3.75
$ws.Cells.Item(245, 13).Value = 1.6
$ws.Cells.Item(245, 14).Value = 6
$ws.Cells.Item(245, 15).Value = 3.75
$ws.Cells.Item(245, 16).Value = 1.571
$ws.Cells.Item(245, 17).Value = 1
$ws.Cells.Item(245, 18).Value = 1.8
$ws.Cells.Item(245, 19).Value = 2.05
$ws.Cells.Item(245, 20).Value = 2.25
$ws.Cells.Item(245, 21).Value = 1.925
$ws.Cells.Item(245, 22).Value = 1.925
$ws.Cells.Item(245, 23).Value = 0
$ws.Cells.Item(245, 24).Value = 0
$ws.Cells.Item(245, 25).Value = 0
$ws.Cells.Item(245, 26).Value = 0
$ws.Cells.Item(245, 27).Value = 0

# Row 246
$ws.Range("A240:AC240").Copy() | Out-Null
$ws.Range("A246:AC246").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(246, 1).Value = 244
$ws.Cells.Item(246, 2).Value = 6978462
$ws.Cells.Item(246, 3).Value = "Bulgaria First League"
$ws.Cells.Item(246, 4).Value = "Bulgaria First League"
$ws.Cells.Item(246, 5).Value = 45402.59375
$ws.Cells.Item(246, 6).Value = "Levski Sofia"
$ws.Cells.Item(246, 7).Value = "Beroe"
$ws.Cells.Item(246, 11).Value = 1.444
$ws.Cells.Item(246, 12).Value = 4.2
$ws.Cells.Item(246, 13).Value = 7.5
$ws.Cells.Item(246, 14).Value = 1.333
$ws.Cells.Item(246, 15).Value = 4.5
$ws.Cells.Item(246, 16).Value = 8.5
$ws.Cells.Item(246, 17).Value = -1.25
$ws.Cells.Item(246, 18).Value = 1.825
$ws.Cells.Item(246, 19).Value = 2.025
$ws.Cells.Item(246, 20).Value = 2.25
$ws.Cells.Item(246, 21).Value = 1.8
$ws.Cells.Item(246, 22).Value = 2.05
$ws.Cells.Item(246, 23).Value = 0
$ws.Cells.Item(246, 24).Value = 0
$ws.Cells.Item(246, 25).Value = 0
$ws.Cells.Item(246, 26).Value = 0
$ws.Cells.Item(246, 27).Value = 0

# Row 247
$ws.Range("A240:AC240").Copy() | Out-Null
$ws.Range("A247:AC247").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(247, 1).Value = 245
$ws.Cells.Item(247, 2).Value = 6978463
$ws.Cells.Item(247, 3).Value = "Bulgaria First League"
$ws.Cells.Item(247, 4).Value = "Bulgaria First League"
$ws.Cells.Item(247, 5).Value = 45403.38541666666
$ws.Cells.Item(247, 6).Value = "Lokomotiv Plovdiv"
$ws.Cells.Item(247, 7).Value = "CSKA 1948 Sofia"
$ws.Cells.Item(247, 11).Value = 1.909
$ws.Cells.Item(247, 12).Value = 3.4
$ws.Cells.Item(247, 13).Value = 4
$ws.Cells.Item(247, 14).Value = 1.8
$ws.Cells.Item(247, 15).Value = 3.5
$ws.Cells.Item(247, 16).Value = 4.5
$ws.Cells.Item(247, 17).Value = -0.75
$ws.Cells.Item(247, 18).Value = 2.05
$ws.Cells.Item(247, 19).Value = 1.8
$ws.Cells.Item(247, 20).Value = 2.25
$ws.Cells.Item(247, 21).Value = 1.85
$ws.Cells.Item(247, 22).Value = 2
$ws.Cells.Item(247, 23).Value = 0
$ws.Cells.Item(247, 24).Value = 0
$ws.Cells.Item(247, 25).Value = 0
$ws.Cells.Item(247, 26).Value = 0
$ws.Cells.Item(247, 27).Value = 0

$excel.CutCopyMode = 0